# Actualizar precios con datos nuevos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the timestamp on the existing row (re-saved capture, same instant to the second).
$ws.Range("A2").Value = 45804.53517370371

# New price-check row: EVOWHEY PROTEIN, 2Kg, 37,90€, captured earlier the same day.
$ws.Range("A3").Value = 45804.45369941551
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B3").Value = "EVOWHEY PROTEIN"
$ws.Range("C3").Value = "2Kg"
$ws.Range("D3").Value = "37,90€"
